# Apply the "update schedule, add homework guidelines" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab to match the workbook/file name.
$ws.Name = "class_schedule_xlsx"

# 2. Update the "Do Before Class" cell for the Pandas Cleaning / Tracebacks
#    row (row 10) to add the new Homework Guidelines + Tracebacks links,
#    keeping the rest of the existing bullet list intact (with "Editing
#    Values" moved to the end). A leading apostrophe forces Excel to treat
#    the leading "-" as literal text (quote-prefixed) instead of a formula,
#    matching the original cell's formatting.
$c10Lines = @(
    '- `Homework Guidelines <homework_guidelines.ipynb>`_',
    '- `What are Tracebacks? <https://www.youtube.com/watch?v=JD8BrXXNtjA>`_',
    '- WM Chapter 6',
    '- WM Chapter 7',
    '- `Python Strings (string section only!) <https://realpython.com/python-data-types/#strings>`_',
    '- `Identifying Problems <cleaning_identifying.ipynb>`_',
    '- `Editing Values <cleaning_editingvalues.ipynb>`_'
)
$ws.Range("C10").Value = "'" + ($c10Lines -join "`n")

# 3. Grow row 10 to fit the extra line that was added to the cell above.
$ws.Rows.Item(10).RowHeight = 99

# 4. Move the active selection to C11 (where the author ended up working).
$ws.Range("C11").Select() | Out-Null
